# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" list (years) for worker ONEYDA FRANCISCA ARROYO RUIZ
# is re-sorted from descending (2105..2005) to ascending (2005..2105), and
# the two "Valor Mora" amounts that sit on the first/last row of that block
# (F16 / F28) are swapped along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periods = @("2005","2006","2007","2008","2009","2010","2011","2012","2101","2102","2103","2104","2105")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# The amount that used to belong to the "2105" row (now row 28) moves to the
# "2005" row (now row 16), and vice versa.
$ws.Range("F16").Value = 26919
$ws.Range("F28").Value = 28090
